# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to the Sheets workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 7063.3
$ws.Range("J32").Value = 7516.75
$ws.Range("L32").Value = 7516.75
$ws.Range("N32").Value = -8168.75
$ws.Range("H74").Value = 5799.2666
$ws.Range("I74").Value = 5272.25
$ws.Range("K74").Value = 5272.25
$ws.Range("M74").Value = -4336.25
$ws.Range("H77").Value = 5799.2666
$ws.Range("I77").Value = 5272.25
$ws.Range("K77").Value = 26361.25
$ws.Range("M77").Value = -21681.25
$ws.Range("H80").Value = 1317.5807
$ws.Range("I80").Value = 523.5
$ws.Range("J80").Value = 1971.5294
$ws.Range("K80").Value = 1570.5
$ws.Range("L80").Value = 5914.5882
$ws.Range("M80").Value = -572.5
$ws.Range("N80").Value = -7910.5882
$ws.Range("H83").Value = 1317.5807
$ws.Range("I83").Value = 523.5
$ws.Range("J83").Value = 1971.5294
$ws.Range("K83").Value = 4711.5
$ws.Range("L83").Value = 17743.7646
$ws.Range("M83").Value = 280.5
$ws.Range("N83").Value = -27727.7646
$ws.Range("H86").Value = 4291.0435
$ws.Range("I86").Value = 3635.9285
$ws.Range("K86").Value = 3635.9285
$ws.Range("M86").Value = -2512.9285
$ws.Range("H89").Value = 4291.0435
$ws.Range("I89").Value = 3635.9285
$ws.Range("K89").Value = 18179.6425
$ws.Range("M89").Value = -12563.6425
$ws.Range("H99").Value = 609.0909
$ws.Range("I99").Value = 563.75
$ws.Range("J99").Value = 730
$ws.Range("K99").Value = 1691.25
$ws.Range("L99").Value = 2190
$ws.Range("M99").Value = -193.25
$ws.Range("N99").Value = -5186
$ws.Range("H111").Value = 2821.6428
$ws.Range("I111").Value = 3451.6667
$ws.Range("K111").Value = 10355.0001
$ws.Range("M111").Value = -7288.000100000001
$ws.Range("H132").Value = 38465120
$ws.Range("I132").Value = 43481748
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 130445244
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -130442714
$ws.Range("N132").Value = -18059

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H34").Value = 210000
$ws.Range("J34").Value = 265000
$ws.Range("L34").Value = 265000
$ws.Range("N34").Value = -265542
$ws.Range("H45").Value = 2396.6086
$ws.Range("I45").Value = 2140.1875
$ws.Range("K45").Value = 2140.1875
$ws.Range("M45").Value = -1763.1875
$ws.Range("H61").Value = 55557670
$ws.Range("J61").Value = 3225
$ws.Range("L61").Value = 3225
$ws.Range("N61").Value = -3649
$ws.Range("H74").Value = 30305466
$ws.Range("I74").Value = 111112880
$ws.Range("K74").Value = 111112880
$ws.Range("M74").Value = -111112006
$ws.Range("H77").Value = 30305466
$ws.Range("I77").Value = 111112880
$ws.Range("K77").Value = 555564400
$ws.Range("M77").Value = -555560032
$ws.Range("H136").Value = 55557670
$ws.Range("J136").Value = 3225
$ws.Range("L136").Value = 9675
$ws.Range("N136").Value = -14775

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 8369.736999999999
$ws.Range("I86").Value = 9852.429
$ws.Range("J86").Value = 4218.2
$ws.Range("K86").Value = 9852.429
$ws.Range("L86").Value = 4218.2
$ws.Range("M86").Value = -8729.429
$ws.Range("N86").Value = -6464.2
$ws.Range("H89").Value = 8369.736999999999
$ws.Range("I89").Value = 9852.429
$ws.Range("J89").Value = 4218.2
$ws.Range("K89").Value = 49262.145
$ws.Range("L89").Value = 21091
$ws.Range("M89").Value = -43646.145
$ws.Range("N89").Value = -32323
$ws.Range("H99").Value = 3829.6667
$ws.Range("I99").Value = 3273.111
$ws.Range("K99").Value = 3273.111
$ws.Range("M99").Value = -1775.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 190.82857
$ws.Range("I7").Value = 92
$ws.Range("J7").Value = 249.22728
$ws.Range("K7").Value = 92
$ws.Range("L7").Value = 249.22728
$ws.Range("M7").Value = 21
$ws.Range("N7").Value = -475.22728
$ws.Range("H58").Value = 4618.467
$ws.Range("I58").Value = 4736.273
$ws.Range("K58").Value = 4736.273
$ws.Range("M58").Value = -4533.273
$ws.Range("H134").Value = 2368.1904
$ws.Range("I134").Value = 2236.6
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 6709.799999999999
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -4174.799999999999
$ws.Range("N134").Value = -20070
$ws.Range("H136").Value = 4618.467
$ws.Range("I136").Value = 4736.273
$ws.Range("K136").Value = 14208.819
$ws.Range("M136").Value = -11658.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 150
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 150
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H134").Value = 2182.75
$ws.Range("I134").Value = 849.3570999999999
$ws.Range("K134").Value = 2548.0713
$ws.Range("M134").Value = 2521.9287
$ws.Range("H139").Value = 2013.0454
$ws.Range("I139").Value = 1911.45
$ws.Range("J139").Value = 3029
$ws.Range("K139").Value = 5734.35
$ws.Range("L139").Value = 9087
$ws.Range("M139").Value = -594.3500000000004
$ws.Range("N139").Value = -19367
$ws.Range("H140").Value = 1473.68
$ws.Range("I140").Value = 1239.8182
$ws.Range("K140").Value = 3719.4546
$ws.Range("M140").Value = 1460.5454
$ws.Range("H141").Value = 9572.923000000001
$ws.Range("I141").Value = 7055.125
$ws.Range("K141").Value = 21165.375
$ws.Range("M141").Value = -15985.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3919.5356
$ws.Range("I70").Value = 3725.5715
$ws.Range("J70").Value = 4113.5
$ws.Range("K70").Value = 3725.5715
$ws.Range("L70").Value = 4113.5
$ws.Range("M70").Value = -3455.5715
$ws.Range("N70").Value = -4653.5
$ws.Range("H73").Value = 3919.5356
$ws.Range("I73").Value = 3725.5715
$ws.Range("J73").Value = 4113.5
$ws.Range("K73").Value = 3725.5715
$ws.Range("L73").Value = 4113.5
$ws.Range("M73").Value = -2789.5715
$ws.Range("N73").Value = -5985.5
$ws.Range("H102").Value = 3769.3225
$ws.Range("I102").Value = 2891.9412
$ws.Range("K102").Value = 2891.9412
$ws.Range("M102").Value = -1269.9412

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 894046.2
$ws.Range("I93").Value = 2818
$ws.Range("K93").Value = 2818
$ws.Range("M93").Value = -1570
$ws.Range("H100").Value = 5056.467
$ws.Range("I100").Value = 3602.5557
$ws.Range("K100").Value = 3602.5557
$ws.Range("M100").Value = -3061.5557
$ws.Range("H122").Value = 5004441
$ws.Range("I122").Value = 3951.8125
$ws.Range("K122").Value = 11855.4375
$ws.Range("M122").Value = -9405.4375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 40002790
$ws.Range("I122").Value = 41667904
$ws.Range("K122").Value = 125003712
$ws.Range("M122").Value = -125001262
$ws.Range("H136").Value = 2349.5715
$ws.Range("I136").Value = 2241.1177
$ws.Range("K136").Value = 6723.353099999999
$ws.Range("M136").Value = -4173.353099999999
